# Ensure the workbook / application objects are freshly (re)bound.
$wb = $excel.ActiveWorkbook

# Add the new worksheet "pagina5" as the last (6th) sheet, right after
# the current last sheet ("pagina4").
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add($null, $lastSheet)
$newSheet.Name = "pagina5"

# Header row (row 1): Document / Description / URL in columns B:D.
$newSheet.Range("B1").Value = "Document"
$newSheet.Range("C1").Value = "Description"
$newSheet.Range("D1").Value = "URL"

$headerRange = $newSheet.Range("B1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$data = @(
    @('MIL-PRF-19500/426J w/Amendment 2 (Initial Draft) -- Dated 3/31/2023 ', 'Transistor, PNP, Silicon, Amplifier Type 2N4957, JAN, JANTX, JANTXV, JANS, JANHC, JANKC ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idprf19500ss426.pdf'),
    @('MIL-PRF-19500/439J w/Amendment 1 (Initial Draft) -- Dated 4/3/2023 ', 'Semiconductor Device, Transistor, NPN, Silicon, High-Power, Types 2N5038 and 2N5039, JAN, JANTX, JANTXV, JANS, JANHC, and JANKC ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idprf19500ss439.pdf'),
    @('MIL-PRF-19500/782 (Initial Draft) -- Dated 4/3/2023 ', 'Transistor, Gallium Nitride, High Electron Mobility Transistor (HEMT), Radiation Hardened, Enhancement Mode, Types 2N7667UFB, 2N7668UFB, 2N7669UFB, Quality Levels JANTXV, JANS JANHC, and JANKC ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-19500/idprf19500ss782.pdf'),
    @('MIL-PRF-39016/48E -- Dated 3/8/2023 ', 'Relays, Electromagnetic, Established Reliability, DPDT, Low Level to 0.5 Ampere (.100 D.I.P. Terminal Spacing), One-Tenth Size, Sensitive, Monostable ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-39016/prf39016ss48.pdf'),
    @('MIL-PRF-55339/39A w/Amendment 1 -- Dated 2/28/2023 ', 'Adapter, Connector, Coaxial, Radio Frequency, (Between Series BNC to Series TNC), Class 2, Straight Plug ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf55339ss39.pdf'),
    @('MIL-PRF-55339/48B w/Amendment 2 -- Dated 2/28/2023 ', 'Adapter, Connector, Coaxial, Radio Frequency, (Between Series SMA to Series TNC), Class 2, Straight Plug ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf55339ss48.pdf'),
    @('MIL-PRF-55339/51A w/Amendment 1 -- Dated 2/28/2023 ', 'Adapter, Connector, Coaxial, Radio Frequency, (Between Series TNC to Series N), Class 2, Straight Plug ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf55339ss51.pdf'),
    @('MIL-PRF-55339/54B w/Amendment 2 -- Dated 2/1/2023 ', 'Adapter, Connector, Electrical, Coaxial, Radio Frequency, (Between Series SMA to N) ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf55339ss54.pdf'),
    @('MIL-PRF-55339/55 w/Amendment 2 -- Dated 2/28/2023 ', 'Adapter, Connectors, Electrical, Coaxial, Radio Frequency, Series SMA, Connector Saver ', 'https://landandmaritimeapps.dla.mil/Downloads/MilSpec/Docs/MIL-PRF-55339/prf55339ss55.pdf')
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $newSheet.Cells.Item($r, 1).Value = $i
    $newSheet.Cells.Item($r, 2).Value = $data[$i][0]
    $newSheet.Cells.Item($r, 3).Value = $data[$i][1]
    $newSheet.Cells.Item($r, 4).Value = $data[$i][2]

    $aCell = $newSheet.Cells.Item($r, 1)
    $aCell.Font.Bold = $true
    $aCell.HorizontalAlignment = -4108
    $aCell.VerticalAlignment = -4160
    $aCell.Borders.LineStyle = 1
}

Write-Output ("Added sheet " + $newSheet.Name + " with " + $wb.Worksheets.Count + " total worksheets")
